$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2023-08-25 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-26 Saturday", 2) | Out-Null

# Update each arithmetic-problem cell directly by table coordinates
# (avoids accidental partial-text collisions between similar expressions,
#  e.g. "6+35=" being a substring of "56+35=")
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "92-3="  # was "47+27="
$t.Cell(1, 2).Range.Text = "7+89="  # was "76-67="
$t.Cell(1, 3).Range.Text = "94-75="  # was "63+19="
$t.Cell(1, 4).Range.Text = "35+46="  # was "92-63="
$t.Cell(1, 5).Range.Text = "59+27="  # was "18+59="
$t.Cell(2, 1).Range.Text = "36+5="  # was "32-28="
$t.Cell(2, 2).Range.Text = "40-38="  # was "76-59="
$t.Cell(2, 3).Range.Text = "40-11="  # was "50-25="
$t.Cell(2, 4).Range.Text = "56+27="  # was "69+22="
$t.Cell(2, 5).Range.Text = "37+56="  # was "86-67="
$t.Cell(3, 1).Range.Text = "70-47="  # was "67+9="
$t.Cell(3, 2).Range.Text = "45+28="  # was "93-34="
$t.Cell(3, 3).Range.Text = "39+7="  # was "77-69="
$t.Cell(3, 4).Range.Text = "7+15="  # was "63-19="
$t.Cell(3, 5).Range.Text = "34-9="  # was "82-48="
$t.Cell(4, 1).Range.Text = "61-3="  # was "7+87="
$t.Cell(4, 2).Range.Text = "48+7="  # was "42-5="
$t.Cell(4, 3).Range.Text = "80-31="  # was "34+47="
$t.Cell(4, 4).Range.Text = "38+58="  # was "29+23="
$t.Cell(4, 5).Range.Text = "26+27="  # was "28+57="
$t.Cell(5, 1).Range.Text = "90-77="  # was "29+43="
$t.Cell(5, 2).Range.Text = "38+24="  # was "26+46="
$t.Cell(5, 3).Range.Text = "47+17="  # was "16+69="
$t.Cell(5, 4).Range.Text = "44-7="  # was "58+33="
$t.Cell(5, 5).Range.Text = "23+59="  # was "71-65="
$t.Cell(6, 1).Range.Text = "48+3="  # was "7+79="
$t.Cell(6, 2).Range.Text = "80-39="  # was "35+37="
$t.Cell(6, 3).Range.Text = "70-39="  # was "92-29="
$t.Cell(6, 4).Range.Text = "47+14="  # was "63-58="
$t.Cell(6, 5).Range.Text = "60-13="  # was "13+18="
$t.Cell(7, 1).Range.Text = "45+7="  # was "19+66="
$t.Cell(7, 2).Range.Text = "48+26="  # was "90-14="
$t.Cell(7, 3).Range.Text = "29+49="  # was "40-1="
$t.Cell(7, 4).Range.Text = "55+17="  # was "60-57="
$t.Cell(7, 5).Range.Text = "75+8="  # was "39+49="
$t.Cell(8, 1).Range.Text = "33+19="  # was "90-18="
$t.Cell(8, 2).Range.Text = "3+18="  # was "6+35="
$t.Cell(8, 3).Range.Text = "45+17="  # was "79+3="
$t.Cell(8, 4).Range.Text = "52-15="  # was "70-54="
$t.Cell(8, 5).Range.Text = "24+39="  # was "93-75="
$t.Cell(9, 1).Range.Text = "52+29="  # was "92-37="
$t.Cell(9, 2).Range.Text = "71-34="  # was "31-5="
$t.Cell(9, 3).Range.Text = "51-13="  # was "56+35="
$t.Cell(9, 4).Range.Text = "82-4="  # was "60-12="
$t.Cell(9, 5).Range.Text = "70-65="  # was "35+56="
$t.Cell(10, 1).Range.Text = "17+54="  # was "41-36="
$t.Cell(10, 2).Range.Text = "82-69="  # was "81-26="
$t.Cell(10, 3).Range.Text = "88-19="  # was "78-29="
$t.Cell(10, 4).Range.Text = "5+86="  # was "59+39="
$t.Cell(10, 5).Range.Text = "15+78="  # was "35+48="
$t.Cell(11, 1).Range.Text = "27+18="  # was "28+4="
$t.Cell(11, 2).Range.Text = "61-44="  # was "15+69="
$t.Cell(11, 3).Range.Text = "34+8="  # was "70-67="
$t.Cell(11, 4).Range.Text = "28+13="  # was "66-7="
$t.Cell(11, 5).Range.Text = "64-19="  # was "53-34="
$t.Cell(12, 1).Range.Text = "82-33="  # was "78-19="
$t.Cell(12, 2).Range.Text = "92-73="  # was "9+39="
$t.Cell(12, 3).Range.Text = "31-13="  # was "7+59="
$t.Cell(12, 4).Range.Text = "38+3="  # was "42-28="
$t.Cell(12, 5).Range.Text = "70-48="  # was "54-29="
$t.Cell(13, 1).Range.Text = "26+59="  # was "52-27="
$t.Cell(13, 2).Range.Text = "96-48="  # was "93-5="
$t.Cell(13, 3).Range.Text = "17+16="  # was "7+38="
$t.Cell(13, 4).Range.Text = "39+24="  # was "37+39="
$t.Cell(13, 5).Range.Text = "46+49="  # was "74-57="
$t.Cell(14, 1).Range.Text = "28+39="  # was "22+59="
$t.Cell(14, 2).Range.Text = "40-15="  # was "92-45="
$t.Cell(14, 3).Range.Text = "82-79="  # was "24+37="
$t.Cell(14, 4).Range.Text = "56-19="  # was "80-24="
$t.Cell(14, 5).Range.Text = "79+7="  # was "74-56="
$t.Cell(15, 1).Range.Text = "76-49="  # was "8+65="
$t.Cell(15, 2).Range.Text = "69+27="  # was "18+77="
$t.Cell(15, 3).Range.Text = "50-3="  # was "51-45="
$t.Cell(15, 4).Range.Text = "60-24="  # was "72-14="
$t.Cell(15, 5).Range.Text = "61-14="  # was "91-23="
$t.Cell(16, 1).Range.Text = "47+37="  # was "75-48="
$t.Cell(16, 2).Range.Text = "20-16="  # was "63-39="
$t.Cell(16, 3).Range.Text = "13+29="  # was "43-24="
$t.Cell(16, 4).Range.Text = "62-19="  # was "24+49="
$t.Cell(16, 5).Range.Text = "80-61="  # was "8+49="
$t.Cell(17, 1).Range.Text = "95-37="  # was "89+7="
$t.Cell(17, 2).Range.Text = "17+34="  # was "6+15="
$t.Cell(17, 3).Range.Text = "49+7="  # was "69+9="
$t.Cell(17, 4).Range.Text = "74-48="  # was "7+46="
$t.Cell(17, 5).Range.Text = "70-34="  # was "16-8="
$t.Cell(18, 1).Range.Text = "83-76="  # was "38+6="
$t.Cell(18, 2).Range.Text = "84-65="  # was "24+19="
$t.Cell(18, 3).Range.Text = "61-42="  # was "53+39="
$t.Cell(18, 4).Range.Text = "9+27="  # was "61-8="
$t.Cell(18, 5).Range.Text = "80-12="  # was "85-19="
$t.Cell(19, 1).Range.Text = "65-48="  # was "49+15="
$t.Cell(19, 2).Range.Text = "49+8="  # was "91-2="
$t.Cell(19, 3).Range.Text = "36+9="  # was "80-22="
$t.Cell(19, 4).Range.Text = "72-54="  # was "49+9="
$t.Cell(19, 5).Range.Text = "5+89="  # was "25+7="
$t.Cell(20, 1).Range.Text = "59+35="  # was "6+18="
$t.Cell(20, 2).Range.Text = "85-68="  # was "29+9="
$t.Cell(20, 3).Range.Text = "75-47="  # was "6+89="
$t.Cell(20, 4).Range.Text = "63-9="  # was "18+37="
$t.Cell(20, 5).Range.Text = "87-79="  # was "51-25="

Write-Host "Replacements complete"
